# Generate Report for Handback
# -----------------------------------------------------------------------
# This script mirrors a "handback" localization-status report refresh:
#   * Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#     on every sheet that shows it (Overview + per-locale sheets).
#   * Each per-locale sheet (zh-cn, de-de) gets its "Latest Target File" /
#     "Latest Handback File" / "Latest Handback DateTime" columns populated
#     for both data rows, with hyperlinks added on the newly-filled target
#     file cells (matching the existing source-file hyperlinks).
#   * The widened columns (status-ish / target-file columns) get resized to
#     fit the new, longer content.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: status column (E/F) text + wider columns
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

$overview.Columns("E").ColumnWidth = 29.17
$overview.Columns("F").ColumnWidth = 29.17

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

$zhcn.Range("J2").Value = "04f251d7-d6e7-4c82-95c2-72df0bf1295c.cebe8c11050796873f842b5764575a3ad603c74d.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-29 00:29:56"

$zhcn.Range("J3").Value = "a2bb649c-4c75-498c-abd5-b355b49a8c89.3b5dfb981a1d8ad4ec3337c92543197284af7f7b.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-29 00:29:56"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/50e2993f5ac0298519dc1faf32f1346b03707929/e2e/04f251d7-d6e7-4c82-95c2-72df0bf1295c.md", "", "", "04f251d7-d6e7-4c82-95c2-72df0bf1295c.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/50e2993f5ac0298519dc1faf32f1346b03707929/e2e/a2bb649c-4c75-498c-abd5-b355b49a8c89.md", "", "", "a2bb649c-4c75-498c-abd5-b355b49a8c89.md")

$zhcn.Columns("C").ColumnWidth = 29.17
$zhcn.Columns("I").ColumnWidth = 39.17
$zhcn.Columns("J").ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

$dede.Range("J2").Value = "04f251d7-d6e7-4c82-95c2-72df0bf1295c.cebe8c11050796873f842b5764575a3ad603c74d.de-de.xlf"
$dede.Range("K2").Value = "2016-08-29 00:30:11"

$dede.Range("J3").Value = "a2bb649c-4c75-498c-abd5-b355b49a8c89.3b5dfb981a1d8ad4ec3337c92543197284af7f7b.de-de.xlf"
$dede.Range("K3").Value = "2016-08-29 00:30:11"

$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/50e2993f5ac0298519dc1faf32f1346b03707929/e2e/04f251d7-d6e7-4c82-95c2-72df0bf1295c.md", "", "", "04f251d7-d6e7-4c82-95c2-72df0bf1295c.md")
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/50e2993f5ac0298519dc1faf32f1346b03707929/e2e/a2bb649c-4c75-498c-abd5-b355b49a8c89.md", "", "", "a2bb649c-4c75-498c-abd5-b355b49a8c89.md")

$dede.Columns("C").ColumnWidth = 29.17
$dede.Columns("I").ColumnWidth = 39.17
$dede.Columns("J").ColumnWidth = 39.17
